# VCHAINS workbook refresh: re-pull of Value Chains data for company 4296954028
# (previously saved under company id 4295891508 / "Lite-On Technology Corp",
# now reflecting "Japan Display Inc"). Updates the run timestamp, the
# criteria values, the full data table (rows 7-11) and a few column widths
# that shifted because of the new (generally shorter) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Run timestamp / criteria -------------------------------------------------
$ws.Range("B1").Value = 44597.8576388889
$ws.Range("B3").Value = "Japan Display Inc"
$ws.Range("B4").Value = 4296954028

# --- Row 7 (Apple Inc / Customer) - only freshness metrics changed -----------
$ws.Range("H7").Value = 44028
$ws.Range("I7").Value = 569
$ws.Range("J7").Value = 4
$ws.Range("K7").Value = 77

# --- Row 8 (now LG Display Co Ltd / Customer, Korea) --------------------------
$ws.Range("A8").Value = 4295882602
$ws.Range("B8").Value = "LG Display Co Ltd"
$ws.Range("D8").Value = "Customer"
$ws.Range("E8").Value = "Korea; Republic (S. Korea)"
$ws.Range("F8").Value = "Electronic Equipment & Parts"
$ws.Range("G8").Value = 0.30295536
$ws.Range("H8").Value = 43487
$ws.Range("I8").Value = 1110
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 26127175997.2868
$ws.Range("M8").Value = 43
$ws.Range("N8").Value = "BB"

# --- Row 9 (now TPK Holding Co Ltd / Supplier, Taiwan) -------------------------
$ws.Range("A9").Value = 4297787586
$ws.Range("B9").Value = "TPK Holding Co Ltd"
$ws.Range("C9").Value = "Public"
$ws.Range("E9").Value = "Taiwan"
$ws.Range("F9").Value = "Electronic Equipment & Parts"
$ws.Range("G9").Value = 0.301104
$ws.Range("H9").Value = 43658
$ws.Range("I9").Value = 939
$ws.Range("J9").Value = 3
$ws.Range("L9").Value = 3893466268.34875
$ws.Range("M9").Value = 41
$ws.Range("N9").Value = "A-"

# --- Row 10 (now IQE PLC / Customer, United Kingdom) ---------------------------
$ws.Range("A10").Value = 4295896478
$ws.Range("B10").Value = "IQE PLC"
$ws.Range("C10").Value = "Public"
$ws.Range("D10").Value = "Customer"
$ws.Range("E10").Value = "United Kingdom"
$ws.Range("F10").Value = "Semiconductors"
$ws.Range("G10").Value = 0.28509432
$ws.Range("H10").Value = 43418
$ws.Range("I10").Value = 1179
$ws.Range("J10").Value = 3
$ws.Range("L10").Value = 228475432.879399
$ws.Range("M10").Value = 44
$ws.Range("N10").Value = "BB-"

# --- Row 11 (now Lumentum Holdings Inc / Customer, USA) -----------------------
$ws.Range("A11").Value = 5045880046
$ws.Range("B11").Value = "Lumentum Holdings Inc"
$ws.Range("C11").Value = "Public"
$ws.Range("E11").Value = "United States of America"
$ws.Range("F11").Value = "Communications & Networking"
$ws.Range("G11").Value = 0.28391264
$ws.Range("H11").Value = 43418
$ws.Range("I11").Value = 1179
$ws.Range("J11").Value = 3
$ws.Range("L11").Value = 1742800000
$ws.Range("M11").Value = 91
$ws.Range("N11").Value = "BBB+"

# --- Column width tweaks (B, C, F, L got narrower/wider for the new text) -----
$ws.Columns.Item(2).ColumnWidth = 21.833333333333332
$ws.Columns.Item(3).ColumnWidth = 5.166666666666667
$ws.Columns.Item(6).ColumnWidth = 27.666666666666668
$ws.Columns.Item(12).ColumnWidth = 17.5
